$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated notebook, reran simulation:
#  - two new reflections ("Holden", "Rizzie Spiral") were inserted into the
#    lookup list right after "Spiral5", which bumps every later label down
#    by two positions in the B column of the results table
#  - "Thomas Hex" was renamed to "Matthies Hex"
#  - the simulation now also reports results for "Michael-CCHex" and
#    "Michael-SNHex" (two more rows appended at the bottom)

$labels = @(
    "HKL",
    "Spiral5",
    "Holden",
    "Rizzie Spiral",
    "RotRing OmegaMax-90",
    "Equal Angle",
    "Tilt Rotate",
    "CLR",
    "Rizzie Hex",
    "Matthies Hex",
    "Tilt Rotate_Partial",
    "RotRing OmegaMax-60",
    "Equal Angle_Partial",
    "Rizzie Hex_Partial",
    "ND Single",
    "RD Single",
    "TD Single",
    "Morris Single",
    "Ring Perpendicular to ND",
    "Ring Perpendicular to RD",
    "Ring Perpendicular to TD",
    "OffsetFTD",
    "OffsetATD",
    "OffsetF45",
    "OffsetA45",
    "OffsetFRD",
    "OffsetARD",
    "Gaussian Quadrature",
    "Michael-CCHex",
    "Michael-SNHex"
)

for ($i = 1; $i -lt $labels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $labels[$i]
    for ($col = 3; $col -le 20; $col++) {
        $ws.Cells.Item($row, $col).Value = 1
    }
}

# The two brand-new rows (30 & 31) need the same "index" column formatting
# (bold, bordered, centered) that every other A-column cell already has.
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)
